# Replace the division expressions in the document's table cells.
# Each "old" expression is unique in the document, so a simple
# Find/Replace (MatchWholeWord not needed, but we keep MatchCase on)
# over the whole document content is safe.

$d = $word.ActiveDocument

$replacements = @(
    @("527÷4=", "880÷6="),
    @("432÷6=", "847÷8="),
    @("309÷8=", "219÷5="),
    @("852÷6=", "301÷3="),
    @("290÷4=", "283÷9="),
    @("406÷7=", "664÷6="),
    @("799÷2=", "499÷2="),
    @("403÷5=", "917÷6="),
    @("963÷7=", "830÷6="),
    @("277÷6=", "802÷2="),
    @("430÷6=", "418÷5="),
    @("131÷6=", "376÷7="),
    @("881÷6=", "280÷2="),
    @("347÷6=", "906÷5="),
    @("247÷7=", "904÷4="),
    @("634÷8=", "308÷2="),
    @("418÷4=", "737÷9="),
    @("112÷9=", "897÷6="),
    @("585÷9=", "883÷9="),
    @("524÷5=", "313÷3="),
    @("899÷9=", "991÷5="),
    @("568÷3=", "441÷6="),
    @("937÷7=", "747÷9="),
    @("523÷5=", "237÷8="),
    @("275÷3=", "734÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
